# API: fix account role when register
# Rename the header row from the product/inventory columns to the
# user-registration columns, and drop the now-unused 7th column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "User Code"
$ws.Range("B1").Value = "Fullname"
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "Phone Number"
$ws.Range("E1").Value = "Password"
$ws.Range("F1").Value = "Role"

# Remove column G entirely (was "Quantity") so the used range shrinks
# back down to A1:F1.
$ws.Range("G1").EntireColumn.Delete()
